# Generate Report for Handback
# Adds a second handed-back file (b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce.md) to the
# handback-status workbook, alongside the existing file (whose uuid is renamed
# from c25b457e-... to 39b05a38-...), across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$uuid1 = "39b05a38-cb78-449e-9597-5c512e044c2c"
$uuid2 = "b73dfcf0-6ec0-4e0b-9b79-830ad2e0e8ce"

$hash1 = "4e4aa8e48aa9029afbc2e7c492c2977b72fb3b8d"
$hash2 = "27c19c26c70fb024e3fd9a90fcffe9787a6ada25"

$commit1 = "f1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0"
$commit2 = "a0b1c2d3e4f5a6b7c8d9e0f1a2b3c4d5e6f7a8b9"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")

# Rename the existing file's uuid (row 2) and bump its generate-date.
$wsOverview.Range("A2").Value = "$uuid1.md"

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = "e2e\$uuid1.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit1/e2e/$uuid1.md", "", "", "e2e\$uuid1.md")

$wsOverview.Range("G2").Value = "2016-09-06 07:12:37"

# New row for the second handed-back file.
$wsOverview.Range("A3").Value = "$uuid2.md"

$wsOverview.Range("B3").Value = "e2e\$uuid2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit2/e2e/$uuid2.md", "", "", "e2e\$uuid2.md")

$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-06 07:12:37"

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Sheets.Item("zh-cn")

# Update existing row 2 (renamed uuid + refreshed timestamps).
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("I2").Hyperlinks.Delete()

$wsZh.Range("A2").Value = "$uuid1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")

$wsZh.Range("G2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-06 07:12:32"

$wsZh.Range("I2").Value = "$uuid1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")

$wsZh.Range("J2").Value = "$uuid1.$hash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 07:12:53"

# New row 3 for the second handed-back file.
$wsZh.Range("A3").Value = "$uuid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$uuid2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-06 07:12:32"

$wsZh.Range("I3").Value = "$uuid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")

$wsZh.Range("J3").Value = "$uuid2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 07:12:53"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.ListObjects.Item(1).Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Sheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("I2").Hyperlinks.Delete()

$wsDe.Range("A2").Value = "$uuid1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")

$wsDe.Range("G2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-06 07:12:37"

$wsDe.Range("I2").Value = "$uuid1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commit1/e2e/$uuid1.md", "", "", "$uuid1.md")

$wsDe.Range("J2").Value = "$uuid1.$hash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 07:13:01"

# New row 3 for the second handed-back file.
$wsDe.Range("A3").Value = "$uuid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$uuid2.$hash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-06 07:12:37"

$wsDe.Range("I3").Value = "$uuid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commit2/e2e/$uuid2.md", "", "", "$uuid2.md")

$wsDe.Range("J3").Value = "$uuid2.$hash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 07:13:01"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.ListObjects.Item(1).Resize($wsDe.Range("A1:P3"))
